# Update Data by bot, scripted by HH
# Re-points the single data row (row 2) from the 2020 Q3 report to the
# 2019 Q3 report: new report date, EPS, revenue/profit totals, gross
# margin, several now-missing YoY/BPS/cashflow metrics, and the
# Q/year/timestamp bookkeeping columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# REPORTDATE: 2020-09-30 -> 2019-09-30
$ws.Range("H2").Value = "2019-09-30 00:00:00"

# BASIC_EPS
$ws.Range("I2").Value = 0.3

# TOTAL_OPERATE_INCOME
$ws.Range("K2").Value = 180771449.92

# PARENT_NETPROFIT
$ws.Range("L2").Value = 18509728.09

# YSTZ, SJLTZ, BPS, MGJYXJJE no longer reported for this period -> blank
# (apply a text format first so the cleared cell is retained as an empty
# cell rather than being dropped entirely, then restore the default style)
$ws.Range("N2").NumberFormat = "@"
$ws.Range("N2").Value = ""
$ws.Range("N2").Style = "Normal"
$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = ""
$ws.Range("O2").Style = "Normal"
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = ""
$ws.Range("P2").Style = "Normal"
$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = ""
$ws.Range("Q2").Style = "Normal"

# XSMLL
$ws.Range("R2").Value = 27.9326711615

# ISNEW: 1 -> 0 (keep as text, not a number)
$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = "0"
$ws.Range("AB2").Style = "Normal"

# QDATE
$ws.Range("AC2").Value = "2019Q3"

# DATATYPE
$ws.Range("AD2").Value = "2019年 三季报"

# DATAYEAR (keep as text, not a number)
$ws.Range("AE2").NumberFormat = "@"
$ws.Range("AE2").Value = "2019"
$ws.Range("AE2").Style = "Normal"

# EITIME
$ws.Range("AG2").Value = "2020-12-22 16:06:32"
